# Update the "Organization website" value (B10): www.stat.kg -> www.stat.gov.kg
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "www.stat.gov.kg"

# Re-apply formatting on B4 (indicator name cell) - this mirrors the font touch-up
# recorded in the saved file (a duplicate font entry gets created for the cell).
$rng = $ws.Range("B4")
$rng.WrapText = $true
$rng.VerticalAlignment = -4160
$rng.Font.Name = "Calibri"

# Update the selected / active cell shown when the workbook is reopened.
$ws.Range("B4").Select()

Write-Host "Edits applied"
